$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (translate to generic field names)
$ws.Range("C2").Value = "code_subject"
$ws.Range("D2").Value = "name_subject"
$ws.Range("E2").Value = "credit"

# Update credit value for row 4 (INT4003 / Hệ quản trị cơ sở dữ liệu)
$ws.Range("E4").Value = 2
